$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 29, shifting existing rows 29:84 down to 30:85
$ws.Rows.Item(29).Insert()

# Populate the new row 29 with the new data entry
$ws.Range("A29").Value = 2
$ws.Range("B29").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C29").Value = 'Coquimbo'
$ws.Range("D29").Value = 44665
$ws.Range("E29").Value = 4
$ws.Range("F29").Value = 100112030
$ws.Range("G29").Value = 'Poroto granado'
$ws.Range("H29").Value = 'Sin especificar'
$ws.Range("I29").Value = 'Primera'
$ws.Range("J29").Value = 200
$ws.Range("K29").Value = 20000
$ws.Range("L29").Value = 22000
$ws.Range("M29").Value = 21000
$ws.Range("N29").Value = '$/malla 25 kilos'
$ws.Range("O29").Value = 'Provincia de Limarí'
$ws.Range("P29").Value = 840
$ws.Range("Q29").Value = 25
$ws.Range("R29").Value = 'Hortaliza'

# Ensure the date column keeps the same date/time number format as the other date cells
$ws.Range("D29").NumberFormat = $ws.Range("D30").NumberFormat
